$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with new contact test data (replacing old Abhishek Kumar / Adidas data)
$ws.Range("A2").Value = "Reshma"
$ws.Range("B2").Value = "Khan"
$ws.Range("C2").Value = "Naggaro"
$ws.Range("D2").Value = "Customer"

# Update the selected/active cell to C2 as recorded in the saved view state
$ws.Range("C2").Select()
